$wb = $excel.ActiveWorkbook

# Generate the handback report: update the "Correspond Handoff Datetime" and
# "Correspond Handback DateTime" for the second data row (the af843c2a... file)
# on both the zh-cn and de-de sheets, giving it its own timestamps instead of
# sharing the ones from the first data row.

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "2016-03-19 22:36:59"
$wsZhCn.Range("H3").Value = "2016-03-19 22:37:18"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E3").Value = "2016-03-19 22:37:02"
$wsDeDe.Range("H3").Value = "2016-03-19 22:37:23"
